# Applies the "Rest API Get and Post update" change to the comparison sheet:
#  - Header E1: discountpercentage -> discountPercentage
#  - Rows 2-3 (master/test for product id 2 "iPhone X") replaced with product id 1 "iPhone 9"
#  - Rows 4-5 (master/test for product id 3 "Samsung Universe 9") replaced with product id 5 "Huawei P30"
#  - Two new rows 6-7 (master/test for product id 15 "Eau De Perfume Spray") appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    # Force a purely numeric-looking string to stay stored as TEXT (t="s"),
    # matching how every column in this sheet is stored as a string cell.
    if ($Text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $Range.NumberFormat = "@"
        $Range.Value = $Text
        $Range.NumberFormat = "general"
    }
    else {
        $Range.Value = $Text
    }
}

function Set-RowValues {
    param($ws, $Row, $Values)
    $cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        Set-TextValue $ws.Range($cols[$i] + $Row) $Values[$i]
    }
}

# ---- Header fix --------------------------------------------------------
Set-TextValue $ws.Range("E1") "discountPercentage"

# ---- Row 2 / Row 3 : iPhone 9 (master / test) --------------------------
$iphone = @(
    "1",
    "iPhone 9",
    "An apple mobile which is nothing like apple",
    "549",
    "12.96",
    "4.69",
    "94",
    "Apple",
    "smartphones",
    "https://i.dummyjson.com/data/products/1/thumbnail.jpg",
    "https://i.dummyjson.com/data/products/1/1.jpg, https://i.dummyjson.com/data/products/1/2.jpg, https://i.dummyjson.com/data/products/1/3.jpg, https://i.dummyjson.com/data/products/1/4.jpg, https://i.dummyjson.com/data/products/1/thumbnail.jpg"
)
Set-RowValues $ws 2 $iphone
Set-RowValues $ws 3 $iphone
# row3 (test) differs only in discountPercentage
Set-TextValue $ws.Range("G3") "94.86"

# ---- Row 4 / Row 5 : Huawei P30 (master / test) -------------------------
$huawei = @(
    "5",
    "Huawei P30",
    "Huawei’s re-badged P30 Pro New Edition was officially unveiled yesterday in Germany and now the device has made its way to the UK.",
    "499",
    "10.58",
    "4.09",
    "32",
    "Huawei",
    "smartphones",
    "https://i.dummyjson.com/data/products/5/thumbnail.jpg",
    "https://i.dummyjson.com/data/products/5/1.jpg, https://i.dummyjson.com/data/products/5/2.jpg, https://i.dummyjson.com/data/products/5/3.jpg"
)
Set-RowValues $ws 4 $huawei
Set-RowValues $ws 5 $huawei
# row5 (test) differs only in rating
Set-TextValue $ws.Range("F5") "4.0955"

# ---- Row 6 / Row 7 : Eau De Perfume Spray (master / test, brand new) ----
$perfume = @(
    "15",
    "Eau De Perfume Spray",
    "Genuine  Al-Rehab spray perfume from UAE/Saudi Arabia/Yemen High Quality",
    "30",
    "10.99",
    "4.7",
    "105",
    "Lord - Al-Rehab",
    "fragrances",
    "https://i.dummyjson.com/data/products/15/thumbnail.jpg",
    "https://i.dummyjson.com/data/products/15/1.jpg, https://i.dummyjson.com/data/products/15/2.jpg, https://i.dummyjson.com/data/products/15/3.jpg, https://i.dummyjson.com/data/products/15/4.jpg, https://i.dummyjson.com/data/products/15/thumbnail.jpg"
)
Set-RowValues $ws 6 $perfume
Set-RowValues $ws 7 $perfume
# row7 (test) differs only in description
Set-TextValue $ws.Range("C7") "Genuine  Al-Rehab spray perfume from UAE/Saudi Arabia/Yemen High Qualitydfgd"

# ---- L column (master/test marker + formatting) -------------------------
$ws.Range("L2").Value = "master"
$ws.Range("L3").Value = "test"
$ws.Range("L4").Value = "master"
$ws.Range("L5").Value = "test"

# L6/L7 reuse the master/test look (bold Calibri font over the existing
# alternating fills) already present on L2:L5.
$ws.Range("L2").Copy() | Out-Null
$ws.Range("L6").PasteSpecial(-4122) | Out-Null
$ws.Range("L6").Value = "master"
$ws.Range("L6").Font.Bold = $true

$ws.Range("L3").Copy() | Out-Null
$ws.Range("L7").PasteSpecial(-4122) | Out-Null
$ws.Range("L7").Value = "test"
$ws.Range("L7").Font.Bold = $true

$excel.CutCopyMode = 0
